# Applies the "adate" prompt-type change, plus the CRIANCA_VISIT (model sheet)
# cleanup, as described in the commit: "Changed date format, and mad changes
# in CRIANCA_VISIT".
#
# Summary of edits:
#  1. prompt_types sheet: add a new prompt type row "adate" / string / string /
#     "Save only mm.dd.yyyy with support for ?? at all positions".
#  2. survey sheet: every field that used the "custom_date" (index 32) prompt
#     type now uses the new "adate" prompt type.
#  3. model sheet: every row that used "custom_date" as its type now uses
#     "adate" as well.
#  4. Update the active-cell selections on the affected sheets to match where
#     the author was last working.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. prompt_types: register the new "adate" prompt type (row 5)
# ---------------------------------------------------------------------------
$promptTypes = $wb.Worksheets.Item("prompt_types")
$promptTypes.Cells.Item(5, 1).Value = "adate"
$promptTypes.Cells.Item(5, 2).Value = "string"
$promptTypes.Cells.Item(5, 3).Value = "string"
$promptTypes.Cells.Item(5, 4).Value = "Save only mm.dd.yyyy with support for ?? at all positions"

$promptTypes.Activate()
[void]$promptTypes.Range("C11").Select()

# ---------------------------------------------------------------------------
# 2. survey: change every "custom_date" field (column D) to "adate"
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$surveyRows = @(17,28,48,59,79,90,110,121,141,152,172,183,203,214,234,245,265,276,296,307,327,338,358,369,389,400,420,431,451,462,482,493,513,524,544,555,575,586,606,617)
foreach ($r in $surveyRows) {
    $survey.Cells.Item($r, 4).Value = "adate"
}

$survey.Activate()
[void]$survey.Range("D632").Select()

# ---------------------------------------------------------------------------
# 3. model: change every "custom_date" type (column B) to "adate"
# ---------------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")
$modelRows = @(42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,163,167,171,175,179,183,187,191,195,199,203,207,211,215,219,227,233)
foreach ($r in $modelRows) {
    $model.Cells.Item($r, 2).Value = "adate"
}

$model.Activate()
[void]$model.Range("C239").Select()
